$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Insert 5 new rows right after the header row (row 1) ---
# This pushes the existing 20 data rows (old rows 2-21) down to rows 7-26.
$ws.Range("A2:A6").EntireRow.Insert()
# The inserted rows inherit the header's formatting; strip it back to the
# plain (unstyled) look used by the rest of the data rows.
$ws.Range("A2:H6").ClearFormats()

# --- Step 2: Fill the 5 newly inserted rows (new rows 2-6) with fresh sensor data ---
$topNewRows = @(
    @(0,"falling",-2.143470001220703,2.89409122467041,0.02666953206062343,0.0542142912745475,-0.0145080499351024,0.0216857157647609),
    @(100,"falling",-2.163377571105957,2.95642032623291,0.1989836648106574,-0.0039706239476799,-0.0080939643085002,0.0329867228865623),
    @(200,"falling",-2.130269622802735,2.900108051300049,0.1112725704908369,-0.0224492978304624,0.0106901414692401,0.0164933614432811),
    @(300,"falling",-2.321649217605591,2.889909553527832,-0.01538913398981086,-0.0187841057777404,-0.0056505035609006,0.0255036242306232),
    @(400,"falling",-2.246681690216064,2.861854076385498,0.0640860199928283,-0.008552113547921099,-0.001527163083665,0.0328340083360672)
)

for ($i = 0; $i -lt $topNewRows.Count; $i++) {
    $r = 2 + $i
    $row = $topNewRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}

# --- Step 3: Re-sequence the timestamp column of the (now shifted) previously
#     existing rows so the 100ms cadence continues unbroken (rows 7-26) ---
for ($r = 7; $r -le 26; $r++) {
    $ws.Cells.Item($r, 1).Value = ($r - 2) * 100
}

# --- Step 4: Append 5 brand-new rows of sensor data after the existing data (rows 27-31) ---
$bottomNewRows = @(
    @(2500,"falling",0.7057792663574191,3.61713025569916,0.1973931401968027,0.0029016099870204,0.0401643887162208,0.0038179077673703),
    @(2600,"falling",0.6601259231567385,3.235907554626462,0.04604268670081965,0.0178678091615438,0.0429132841527462,0.027030786499381),
    @(2700,"falling",0.9756811141967799,3.334283685684206,0.05510960519313853,-0.0041233403608202,-0.0798706337809562,-0.0476474873721599),
    @(2800,"falling",0.6400001525878865,3.297549438476561,0.02418547868728584,-0.0029016099870204,0.0207694191485643,0.0029016099870204),
    @(2900,"falling",0.471990585327149,3.205180048942565,-0.01620917022228236,-0.007177666760981,-0.0045814891345798,0.0218384321779012)
)

for ($i = 0; $i -lt $bottomNewRows.Count; $i++) {
    $r = 27 + $i
    $row = $bottomNewRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}
